# "Generate Report for Handoff"
# Adds a new handed-off file (269764a5-...) as row 3 on all three sheets:
#   Overview (sheet1), zh-cn (sheet2), de-de (sheet3)

$wb = $excel.ActiveWorkbook

# ---- literal strings (kept in variables to avoid repetition / typos) ----
$mdFile    = '269764a5-1046-4040-acfd-1454d01e2d2dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$e2eMdFile = 'e2e\269764a5-1046-4040-acfd-1454d01e2d2dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$zhcnXlf   = '269764a5-1046-4040-acfd-1454d01e2d2doooooooooooooooooooooooooooooooooooooooo.353b76c7df7a569b32a4608a2aa8063a7455eff2.zh-cn.xlf'
$dedeXlf   = '269764a5-1046-4040-acfd-1454d01e2d2doooooooooooooooooooooooooooooooooooooooo.353b76c7df7a569b32a4608a2aa8063a7455eff2.de-de.xlf'
$oldDedeXlf = '297bf66f-a4d2-4f0d-b70e-b610560303ecoooooooooooooooooooooooooooooooooooooooo.dde86b9c277e77c0f1c3bde04cf36c0dd547b59e.de-de.xlf'

$statusReady = 'Ready for handoff'
$dateOverview = '2016-09-05 04:32:13'
$dateZhcn     = '2016-09-05 04:32:08'
$dateDede     = '2016-09-05 04:32:13'

$dateFmt = 'yyyy-mm-dd HH:mm:ss'

$hyperlinkBase = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57d7d292d1b18e5b470dd300bbab2bcdd648a7ed/e2e/'

# ================= Overview sheet =================
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $dateOverview
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($hyperlinkBase + $mdFile), [Type]::Missing, [Type]::Missing, $e2eMdFile) | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 17.22
$wsOverview.Columns.Item(6).ColumnWidth = 17.22

# ================= zh-cn sheet =================
$wsZhcn = $wb.Worksheets.Item("zh-cn")
$tblZhcn = $wsZhcn.ListObjects.Item(1)
$tblZhcn.ListRows.Add() | Out-Null

$wsZhcn.Range("B3").Value = ".md"
$wsZhcn.Range("C3").Value = $statusReady
$wsZhcn.Range("D3").Value = "e2e"
$wsZhcn.Range("E3").Value = "ht"
$wsZhcn.Range("F3").Value = "'False"
$wsZhcn.Range("G3").Value = $zhcnXlf
$wsZhcn.Range("H3").Value = $dateZhcn
$wsZhcn.Range("H3").NumberFormat = $dateFmt
$wsZhcn.Range("I3").Value = ""
$wsZhcn.Range("J3").Value = ""
$wsZhcn.Range("K3").Value = $oldDedeXlf
$wsZhcn.Range("K3").NumberFormat = $dateFmt
$wsZhcn.Range("L3").Value = ""
$wsZhcn.Range("M3").Value = $dedeXlf
$wsZhcn.Range("N3").Value = ""
$wsZhcn.Range("O3").Value = "'False"
$wsZhcn.Range("P3").Value = ""

$wsZhcn.Hyperlinks.Add($wsZhcn.Range("A3"), ($hyperlinkBase + $mdFile), [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null

$wsZhcn.Columns.Item(3).ColumnWidth = 17.22

# ================= de-de sheet =================
$wsDede = $wb.Worksheets.Item("de-de")
$tblDede = $wsDede.ListObjects.Item(1)
$tblDede.ListRows.Add() | Out-Null

$wsDede.Range("B3").Value = ".md"
$wsDede.Range("C3").Value = $statusReady
$wsDede.Range("D3").Value = "e2e"
$wsDede.Range("E3").Value = "ht"
$wsDede.Range("F3").Value = "'False"
$wsDede.Range("G3").Value = $dedeXlf
$wsDede.Range("H3").Value = $dateDede
$wsDede.Range("H3").NumberFormat = $dateFmt
$wsDede.Range("I3").Value = ""
$wsDede.Range("J3").Value = ""
$wsDede.Range("K3").Value = $oldDedeXlf
$wsDede.Range("K3").NumberFormat = $dateFmt
$wsDede.Range("L3").Value = ""
$wsDede.Range("M3").Value = $dedeXlf
$wsDede.Range("N3").Value = ""
$wsDede.Range("O3").Value = "'False"
$wsDede.Range("P3").Value = ""

$wsDede.Hyperlinks.Add($wsDede.Range("A3"), ($hyperlinkBase + $mdFile), [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null

$wsDede.Columns.Item(3).ColumnWidth = 17.22
